# Updates the cryptos list sheet:
#  - refresh Price (col D) and Volume(1h) (col E) figures for the existing
#    coin rows (2-46)
#  - insert a new "BabyDogeCoin" row at row 47, which pushes Aptos,
#    TheSandbox, RenderToken and EnergySwap down by one row and drops the
#    previously-last row (Algorand) off the bottom of the 50-row table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Price / Volume(1h) refresh for rows 2-46 -----------------------
# Each entry is RowNumber = @(NewPrice, NewVolume). A $null element means
# that particular column is not touched for that row.
$priceVolumeUpdates = @{
    2 = @("29.406.44", "  +0.19%  ")
    3 = @("1.849.55", "  +0.27%  ")
    4 = @("0.9999", "  +0.10%  ")
    5 = @("240.75", "  +0.24%  ")
    6 = @("0.6304", "  -0.11%  ")
    7 = @($null, "  +0.02%  ")
    8 = @("0.07710", "  +2.34%  ")
    9 = @("0.2940", "  -0.56%  ")
    10 = @("24.51", "  +0.39%  ")
    11 = @("0.07748", "  +0.49%  ")
    12 = @("1.844.36", "  -0.19%  ")
    13 = @("5.026", "  +0.78%  ")
    14 = @("0.00001090", "  +8.83%  ")
    15 = @("0.6804", "  -0.29%  ")
    16 = @("83.65", "  +1.08%  ")
    17 = @("2.102.51", "  -7.15%  ")
    18 = @("6.148", "  +0.45%  ")
    19 = @("29.437.14", "  +0.16%  ")
    20 = @("229.19", "  +0.70%  ")
    21 = @("12.47", "  +0.43%  ")
    22 = @("1.000", "  +0.05%  ")
    23 = @("7.455", "  -1.05%  ")
    24 = @($null, "  +0.02%  ")
    25 = @("157.33", "  +0.02%  ")
    26 = @($null, "  -0.53%  ")
    27 = @("8.359", "  +0.28%  ")
    28 = @($null, $null)
    29 = @("1.315", "  +4.84%  ")
    30 = @($null, "  +0.13%  ")
    31 = @($null, "  +0.16%  ")
    32 = @("4.111", "  -0.30%  ")
    33 = @("4.048", "  +0.78%  ")
    34 = @("1.852", "  +0.50%  ")
    35 = @($null, "  +0.42%  ")
    36 = @("0.7092", "  -0.49%  ")
    37 = @("2.588", "  -0.11%  ")
    38 = @("1.231.53", "  -2.20%  ")
    39 = @($null, "  -0.23%  ")
    40 = @("0.01801", "  -0.37%  ")
    41 = @("6.497", "  +4.46%  ")
    42 = @("0.9138", "  +0.34%  ")
    43 = @($null, "  +0.03%  ")
    44 = @("2.011.38", "  +0.47%  ")
    45 = @("101.52", "  +0.36%  ")
    46 = @("66.27", "  +0.20%  ")
}

foreach ($rowNum in $priceVolumeUpdates.Keys) {
    $newPrice = $priceVolumeUpdates[$rowNum][0]
    $newVolume = $priceVolumeUpdates[$rowNum][1]

    if ($null -ne $newPrice) {
        $priceCell = $ws.Cells.Item($rowNum, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $newPrice
    }

    if ($null -ne $newVolume) {
        $volumeCell = $ws.Cells.Item($rowNum, 5)
        $volumeCell.NumberFormat = "@"
        $volumeCell.Value = $newVolume
    }
}

# --- 2. Insert "BabyDogeCoin" as the new row 47 -------------------------
# This shifts Aptos / TheSandbox / RenderToken / EnergySwap down one row
# and the old Algorand row (previously 51) falls off the bottom of the
# table, so we simply rewrite B/C/D/E for rows 47-51 from the bottom up.

$rows47to51 = @(
    @(47, "BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.00000000121", "  +3.24%  "),
    @(48, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.165", "  +1.82%  "),
    @(49, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4016", "  -0.50%  "),
    @(50, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.688", "  +0.29%  "),
    @(51, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "8.992", "  -1.18%  ")
)

foreach ($entry in $rows47to51) {
    $rowNum = $entry[0]
    $coin = $entry[1]
    $link = $entry[2]
    $price = $entry[3]
    $volume = $entry[4]

    $ws.Cells.Item($rowNum, 2).Value = $coin
    $ws.Cells.Item($rowNum, 3).Value = $link

    $priceCell = $ws.Cells.Item($rowNum, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price

    $volumeCell = $ws.Cells.Item($rowNum, 5)
    $volumeCell.NumberFormat = "@"
    $volumeCell.Value = $volume
}
